$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: A1 reference -> new value. Values that could be
# misread as numbers (and would lose formatting, e.g. "94.60" -> 94.6,
# or "1.49" -> 1.49 as a float dropping the original text form) are
# force-written as text via NumberFormat "@" first.
$updates = [ordered]@{
    'D2' = '43.513.92'
    'E2' = '  -1.37%  '
    'D3' = '2.277.68'
    'E3' = '  +1.05%  '
    'E4' = '  +0.08%  '
    'D5' = '94.47'
    'E5' = '  -4.29%  '
    'D6' = '266.16'
    'E6' = '  -2.34%  '
    'E7' = '  -1.09%  '
    'E8' = '  -0.04%  '
    'D9' = '0.604'
    'E9' = '  -3.59%  '
    'D10' = '44.16'
    'E10' = '  -7.83%  '
    'D11' = '0.0931'
    'E11' = '  -1.42%  '
    'D12' = '7.69'
    'E12' = '  -5.99%  '
    'E13' = '  +0.03%  '
    'D14' = '2.619.02'
    'E14' = '  +1.23%  '
    'D15' = '15.09'
    'E15' = '  -2.38%  '
    'D16' = '0.842'
    'E16' = '  +1.00%  '
    'D17' = '2.283.77'
    'E17' = '  +1.17%  '
    'D18' = '43.465.49'
    'E18' = '  -1.41%  '
    'E19' = '  -0.61%  '
    'D20' = '6.16'
    'E20' = '  -1.15%  '
    'D21' = '72.05'
    'E21' = '  +1.73%  '
    'E22' = '  -1.47%  '
    'D23' = '233.13'
    'E23' = '  -0.60%  '
    'D24' = '8.91'
    'E24' = '  -10.51%  '
    'D25' = '0.999'
    'E25' = '  -0.04%  '
    'E26' = '  -1.00%  '
    'D27' = '11.16'
    'E27' = '  -1.96%  '
    'E28' = '  -1.63%  '
    'E29' = '  -0.13%  '
    'D30' = '39.17'
    'E30' = '  -1.83%  '
    'D31' = '175.40'
    'E31' = '  +0.93%  '
    'D32' = '21.76'
    'E32' = '  +2.50%  '
    'D33' = '0.0878'
    'D34' = '5.30'
    'E34' = '  -6.68%  '
    'E35' = '  +0.00%  '
    'B36' = 'Kaspa'
    'C36' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D36' = '0.107'
    'E36' = '  -5.24%  '
    'B37' = 'VeChain'
    'C37' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D37' = '0.0352'
    'E37' = '  -0.56%  '
    'D38' = '4.39'
    'E38' = '  +0.81%  '
    'E39' = '  -6.92%  '
    'E40' = '  +5.44%  '
    'D41' = '0.233'
    'E41' = '  -6.01%  '
    'E42' = '  +14.59%  '
    'E43' = '  -4.76%  '
    'D44' = '62.79'
    'E44' = '  +1.21%  '
    'D45' = '8.76'
    'E45' = '  +2.73%  '
    'D46' = '5.19'
    'E46' = '  -4.81%  '
    'E47' = '  -1.52%  '
    'D48' = '97.31'
    'E48' = '  -3.17%  '
    'E49' = '  -1.15%  '
    'D50' = '1.48'
    'E50' = '  +3.79%  '
    'D51' = '2.498.37'
    'E51' = '  +0.94%  '
}

$textCells = @(
    'D5'
    'D6'
    'D9'
    'D10'
    'D11'
    'D12'
    'D15'
    'D16'
    'D20'
    'D21'
    'D23'
    'D24'
    'D25'
    'D27'
    'D30'
    'D31'
    'D32'
    'D33'
    'D34'
    'D36'
    'D37'
    'D38'
    'D41'
    'D44'
    'D45'
    'D46'
    'D48'
    'D50'
)

foreach ($ref in $updates.Keys) {
    if ($textCells -contains $ref) {
        $ws.Range($ref).NumberFormat = "@"
    }
    $ws.Range($ref).Value = $updates[$ref]
}
